# Applies the Notes_ApresRat.xlsx edits:
#  - swap the "LastName" / "FirstName" header labels in B1 / C1
#  - shift the CNE identifiers in column A (rows 2-11) down by 10
#  - move the active selection to I11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap B1 and C1 header text (LastName <-> FirstName)
$b1 = $ws.Range("B1").Value2
$c1 = $ws.Range("C1").Value2
$ws.Range("B1").Value = $c1
$ws.Range("C1").Value = $b1

# Decrease each CNE (column A, rows 2-11) by 10
for ($r = 2; $r -le 11; $r++) {
    $current = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $current - 10
}

# Update the active selection to I11
$ws.Range("I11").Select()
